# Update countries & provincias Spain
# Applies the COVID-19 "paises" dataset refresh:
#  - Updated case counters for several countries (rank-ordered by total cases)
#  - A handful of countries changed rank position relative to their neighbours,
#    which is expressed here as writing the correct country name + stats into
#    each affected row (so the row order stays a strictly-sorted list).
#  - The "last updated" timestamp banner in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=17;  Values=@("Austria", 11199, 70, 1749, 9292, 227, 0, 158)},
    @{Row=23;  Values=@("Australia", 5350, 36, 585, 4737, 50, 3, 28)},
    @{Row=58;  Values=@("Ucrania", 942, 45, 19, 900, 0, 1, 23)},
    @{Row=65;  Values=@("Armenia", 736, 73, 43, 686, 30, 0, 7)},
    @{Row=66;  Values=@("Marruecos", 735, 27, 34, 654, 1, 3, 47)},
    @{Row=67;  Values=@("Crucero", 712, 0, 619, 82, 10, 0, 11)},
    @{Row=68;  Values=@("Lituania", 696, 47, 7, 680, 11, 0, 9)},
    @{Row=71;  Values=@("Bosnia y Herzegovina", 543, 10, 20, 507, 4, 0, 16)},
    @{Row=74;  Values=@("Letonia", 493, 35, 31, 462, 3, 0, 0)},
    @{Row=75;  Values=@("Bulgaria", 477, 20, 30, 435, 18, 2, 12)},
    @{Row=105; Values=@("Islas Feroe", 179, 2, 91, 88, 1, 0, 0)},
    @{Row=108; Values=@("Montenegro", 160, 16, 0, 158, 4, 0, 2)},
    @{Row=109; Values=@("Sri Lanka", 151, 0, 22, 125, 5, 0, 4)},
    @{Row=110; Values=@("Georgia", 148, 14, 27, 121, 6, 0, 0)},
    @{Row=111; Values=@("Venezuela", 146, 0, 43, 98, 6, 0, 5)}
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $vals = $r.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($rowNum, $i + 1).Value = $vals[$i]
    }
}

$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 09:20"
